$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 376, which pushes the existing rows
# 376-451 down to become rows 378-453 (preserving all their data/styles).
$ws.Rows("376:377").Insert()

# Populate the two newly inserted rows (376 and 377) with fresh data.

# Row 376
$ws.Range("A376").Value = 9
$ws.Range("B376").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C376").Value = "Metropolitana"
$ws.Range("D376").Value = 44776
$ws.Range("E376").Value = 13
$ws.Range("F376").Value = 100112013
$ws.Range("G376").Value = "Alcachofa"
$ws.Range("H376").Value = "Española"
$ws.Range("I376").Value = "Extra"
$ws.Range("J376").Value = 34
$ws.Range("K376").Value = 18000
$ws.Range("L376").Value = 18000
$ws.Range("M376").Value = 18000
$ws.Range("N376").Value = "$/caja 25 unidades"
$ws.Range("O376").Value = "Provincia del Elquí"
$ws.Range("P376").Value = 18000
$ws.Range("Q376").Value = 1
$ws.Range("R376").Value = "Hortaliza"

# Row 377
$ws.Range("A377").Value = 9
$ws.Range("B377").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C377").Value = "Metropolitana"
$ws.Range("D377").Value = 44776
$ws.Range("E377").Value = 13
$ws.Range("F377").Value = 100112013
$ws.Range("G377").Value = "Alcachofa"
$ws.Range("H377").Value = "Española"
$ws.Range("I377").Value = "Primera"
$ws.Range("J377").Value = 52
$ws.Range("K377").Value = 16000
$ws.Range("L377").Value = 16000
$ws.Range("M377").Value = 16000
$ws.Range("N377").Value = "$/caja 30 unidades"
$ws.Range("O377").Value = "Provincia del Elquí"
$ws.Range("P377").Value = 533
$ws.Range("Q377").Value = 30
$ws.Range("R377").Value = "Hortaliza"
